$p = $ppt.ActivePresentation

# Locate the "Methodology" slide / content placeholder that holds the
# "Pre-processing of the data" bullet and the "Final markdown report" link
# (search by content instead of a hard-coded index, so the script is
# resilient to slide re-ordering).
$targetShape = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $txt = $shape.TextFrame.TextRange.Text
            if ($txt.IndexOf("Final markdown report") -ge 0) {
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# 1) Merge the three runs "Pr" / "e-processing " / "of the data " into a
#    single run of text "Pre-processing of the data " (same visible text,
#    just re-typed as one contiguous run).
$mergedText = "Pre-processing of the data "
$full = $tr.Text
$idx = $full.IndexOf($mergedText)
if ($idx -lt 0) {
    $idx = $full.IndexOf("Pre-processing")
}
$sub = $tr.Characters($idx + 1, $mergedText.Length)
$sub.Text = $mergedText

# 2) Re-point the "Final markdown report (...)" hyperlink from the GitHub
#    .Rmd source link to the rendered rawgit.com HTML report, updating both
#    the displayed URL text and the actual hyperlink target.
$oldUrl = "https://github.com/SnakeAkaPython/FinalReportGeo-ScriptingGRS-51806/blob/master/FinalProjectReport.Rmd"
$newUrl = "https://cdn.rawgit.com/SnakeAkaPython/FinalReportGeo-ScriptingGRS-51806/master/FinalProjectReporttest.html"

$full = $tr.Text
$idx = $full.IndexOf($oldUrl)
$linkSub = $tr.Characters($idx + 1, $oldUrl.Length)
$linkSub.Text = $newUrl

$full = $tr.Text
$idx = $full.IndexOf($newUrl)
$linkSub = $tr.Characters($idx + 1, $newUrl.Length)
$hyperlink = $linkSub.ActionSettings.Item(1).Hyperlink
$hyperlink.Address = $newUrl
